$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3: Rohan Bhambare
$ws.Cells.Item(3,1).Value = "Rohan Bhambare"
$ws.Cells.Item(3,2).Value = "devmailproject0@gmail.com"
$ws.Cells.Item(3,3).Value = "scrypt:32768:8:1`$bJNaf0HoGYnu3OBW`$f88aed9489912cb7c2f14455ee8988f7b64115bbeb798cce209ba979c2b80181a6cceb483b883b7fb68e913b8ea549b4491d3a59b2e9dd9cbbf062554305b8ed"
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 45752.21000189815

# New row 4: Vivek Totre
$ws.Cells.Item(4,1).Value = "Vivek Totre"
$ws.Cells.Item(4,2).Value = "tech@gmail.com"
$ws.Cells.Item(4,3).Value = "scrypt:32768:8:1`$AG7HNQuB46fotiUC`$7613167c7537901ea92074ac3aee512f8334554db730b64e0f3018fa9b5a5568d3c8f192da243c0d5791e044db4a2eaa7476fe5069e1ded61272071a783cf5de"
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 45752.21127630154

# Carry the existing "created_at" date/time style (number format) from E2
# down to the two newly added rows so they render as timestamps too.
$ws.Range("E2").Copy()
$ws.Range("E3:E4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Tiny re-save precision adjustment on the pre-existing timestamp.
$ws.Cells.Item(2,5).Value = 45750.99573267361

# Restore the original selection/active cell.
$ws.Range("A1").Select() | Out-Null
